$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("C2").Value = 47.5
$ws.Range("C3").Value = 56.66666666666666
$ws.Range("C5").Value = 56.66666666666666
$ws.Range("C9").Value = 90
